{"js": "// Apply the textual updates described by the diff: update the date\n// paragraph and every \"A\u00d7B=C\" answer cell in the table, preserving all\n// existing run/paragraph formatting by editing in place via search+replace.\nconst replacements = [\n  [\"2025-10-23 Thursday\", \"2025-10-24 Friday\"],\n  [\"36\u00d723=828\", \"59\u00d718=1062\"],\n  [\"25\u00d793=2325\", \"81\u00d747=3807\"],\n  [\"52\u00d775=3900\", \"37\u00d735=1295\"],\n  [\"83\u00d741=3403\", \"38\u00d782=3116\"],\n  [\"31\u00d736=1116\", \"38\u00d731=1178\"],\n  [\"80\u00d720=1600\", \"83\u00d744=3652\"],\n  [\"46\u00d781=3726\", \"87\u00d714=1218\"],\n  [\"18\u00d728=504\", \"12\u00d722=264\"],\n  [\"18\u00d767=1206\", \"84\u00d716=1344\"],\n  [\"72\u00d773=5256\", \"98\u00d796=9408\"],\n  [\"71\u00d797=6887\", \"71\u00d755=3905\"],\n  [\"92\u00d737=3404\", \"30\u00d729=870\"],\n  [\"29\u00d719=551\", \"65\u00d749=3185\"],\n  [\"60\u00d776=4560\", \"76\u00d727=2052\"],\n  [\"48\u00d753=2544\", \"55\u00d765=3575\"],\n  [\"75\u00d750=3750\", \"72\u00d785=6120\"],\n  [\"73\u00d737=2701\", \"15\u00d791=1365\"],\n  [\"23\u00d730=690\", \"80\u00d798=7840\"],\n  [\"42\u00d764=2688\", \"48\u00d723=1104\"],\n  [\"65\u00d765=4225\", \"29\u00d728=812\"],\n  [\"99\u00d789=8811\", \"68\u00d799=6732\"],\n  [\"84\u00d754=4536\", \"30\u00d782=2460\"],\n  [\"88\u00d737=3256\", \"26\u00d754=1404\"],\n  [\"65\u00d734=2210\", \"47\u00d745=2115\"],\n  [\"25\u00d714=350\", \"81\u00d782=6642\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const rangeItem of results.items) {\n    rangeItem.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the textual updates described by the diff: update the date\n# paragraph and every \"A\u00d7B=C\" answer cell in the table, preserving all\n# existing run/paragraph formatting via Find/Replace over the whole story.\n$d = $word.ActiveDocument\n\n$oldValues = @(\"2025-10-23 Thursday\", \"36\u00d723=828\", \"25\u00d793=2325\", \"52\u00d775=3900\", \"83\u00d741=3403\", \"31\u00d736=1116\", \"80\u00d720=1600\", \"46\u00d781=3726\", \"18\u00d728=504\", \"18\u00d767=1206\", \"72\u00d773=5256\", \"71\u00d797=6887\", \"92\u00d737=3404\", \"29\u00d719=551\", \"60\u00d776=4560\", \"48\u00d753=2544\", \"75\u00d750=3750\", \"73\u00d737=2701\", \"23\u00d730=690\", \"42\u00d764=2688\", \"65\u00d765=4225\", \"99\u00d789=8811\", \"84\u00d754=4536\", \"88\u00d737=3256\", \"65\u00d734=2210\", \"25\u00d714=350\")\n$newValues = @(\"2025-10-24 Friday\", \"59\u00d718=1062\", \"81\u00d747=3807\", \"37\u00d735=1295\", \"38\u00d782=3116\", \"38\u00d731=1178\", \"83\u00d744=3652\", \"87\u00d714=1218\", \"12\u00d722=264\", \"84\u00d716=1344\", \"98\u00d796=9408\", \"71\u00d755=3905\", \"30\u00d729=870\", \"65\u00d749=3185\", \"76\u00d727=2052\", \"55\u00d765=3575\", \"72\u00d785=6120\", \"15\u00d791=1365\", \"80\u00d798=7840\", \"48\u00d723=1104\", \"29\u00d728=812\", \"68\u00d799=6732\", \"30\u00d782=2460\", \"26\u00d754=1404\", \"47\u00d745=2115\", \"81\u00d782=6642\")\n\nfor ($i = 0; $i -lt $oldValues.Count; $i++) {\n    $oldText = $oldValues[$i]\n    $newText = $newValues[$i]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    # wdFindContinue=1, wdReplaceAll=2\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
